# Textbox response formatting fix
$wb = $excel.ActiveWorkbook

# --- Rename worksheets (new randomized timestamp-based task-order IDs) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16511687727284515"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687755140748"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687755150743"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687755776265"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687756495247"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687726861825.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687727111988.csv"
$ws1.Range("B4").Value = "go_stims-1651168772713195.csv"
$ws1.Range("B5").Value = "GNG_stims-16511687727264516.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1651168775499689.csv"
$ws2.Range("B3").Value = "ZB-match_1-1651168773260333.csv"
$ws2.Range("B4").Value = "OB-16511687745661666.csv"
$ws2.Range("B5").Value = "ZB-match_1-16511687729423606.csv"
$ws2.Range("B6").Value = "TB-16511687751260526.csv"
$ws2.Range("B7").Value = "TB-16511687754561496.csv"
$ws2.Range("B8").Value = "OB-16511687744597564.csv"
$ws2.Range("B9").Value = "OB-16511687743108726.csv"
$ws2.Range("B10").Value = "ZB-match_7-16511687730421767.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687755300474.csv"
$ws4.Range("B3").Value = "ZM_stims-16511687755170765.csv"
$ws4.Range("B4").Value = "MM_stims-16511687755610485.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687755300474.csv"
$ws4.Range("B6").Value = "MM_stims-1651168775576623.csv"
$ws4.Range("B7").Value = "ZM_stims-1651168775562062.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16511687755816746.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511687756339335.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687756181574.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687756027293.csv"
